# Auto-generated script to update cryptos list values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.844.32"
$ws.Range("E2").Value = "  +0.12%  "
$ws.Range("D3").Value = "1.637.15"
$ws.Range("E3").Value = "  +0.10%  "
$ws.Range("E4").Value = "  +0.49%  "
$ws.Range("D5").Value = "216.04"
$ws.Range("E5").Value = "  +0.37%  "
$ws.Range("D6").Value = "0.503"
$ws.Range("E6").Value = "  -0.17%  "
$ws.Range("E7").Value = "  +0.46%  "
$ws.Range("E8").Value = "  +0.37%  "
$ws.Range("D9").Value = "0.0638"
$ws.Range("E9").Value = "  -0.45%  "
$ws.Range("D10").Value = "19.65"
$ws.Range("E10").Value = "  -1.90%  "
$ws.Range("D11").Value = "0.0788"
$ws.Range("E11").Value = "  +1.19%  "
$ws.Range("D12").Value = "1.651.50"
$ws.Range("E12").Value = "  +1.02%  "
$ws.Range("E13").Value = "  -0.64%  "
$ws.Range("D14").Value = "1.864.77"
$ws.Range("E14").Value = "  +0.17%  "
$ws.Range("D15").Value = "0.554"
$ws.Range("E15").Value = "  -1.28%  "
$ws.Range("D16").Value = "0.0₃0774"
$ws.Range("E16").Value = "  +1.50%  "
$ws.Range("D17").Value = "63.33"
$ws.Range("E17").Value = "  +0.50%  "
$ws.Range("D18").Value = "25.850.60"
$ws.Range("E18").Value = "  +0.05%  "
$ws.Range("E19").Value = "  +0.43%  "
$ws.Range("E20").Value = "  +2.28%  "
$ws.Range("D21").Value = "193.73"
$ws.Range("E21").Value = "  -0.45%  "
$ws.Range("E22").Value = "  +0.83%  "
$ws.Range("D23").Value = "6.18"
$ws.Range("E23").Value = "  +1.41%  "
$ws.Range("E24").Value = "  +0.57%  "
$ws.Range("D25").Value = "1.78"
$ws.Range("E25").Value = "  +0.16%  "
$ws.Range("D26").Value = "139.95"
$ws.Range("E26").Value = "  -0.25%  "
$ws.Range("D27").Value = "0.120"
$ws.Range("E27").Value = "  -4.22%  "
$ws.Range("D28").Value = "6.84"
$ws.Range("E28").Value = "  +0.58%  "
$ws.Range("D29").Value = "15.61"
$ws.Range("E29").Value = "  +0.95%  "
$ws.Range("E30").Value = "  +0.41%  "
$ws.Range("E31").Value = "  -0.37%  "
$ws.Range("D32").Value = "3.35"
$ws.Range("E32").Value = "  +1.30%  "
$ws.Range("D33").Value = "3.27"
$ws.Range("E33").Value = "  +1.09%  "
$ws.Range("E34").Value = "  +1.36%  "
$ws.Range("E35").Value = "  +0.98%  "
$ws.Range("D36").Value = "0.897"
$ws.Range("D37").Value = "2.59"
$ws.Range("E37").Value = "  +0.52%  "
$ws.Range("D38").Value = "0.551"
$ws.Range("E38").Value = "  -0.33%  "
$ws.Range("D39").Value = "1.109.31"
$ws.Range("E39").Value = "  -1.51%  "
$ws.Range("E40").Value = "  +0.49%  "
$ws.Range("E41").Value = "  +0.62%  "
$ws.Range("D43").Value = "0.807"
$ws.Range("E43").Value = "  +0.92%  "
$ws.Range("D44").Value = "99.57"
$ws.Range("E44").Value = "  +1.59%  "
$ws.Range("E45").Value = "  -2.24%  "
$ws.Range("E46").Value = "  -0.38%  "
$ws.Range("E47").Value = "  +10.04%  "
$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").Value = "0.420"
$ws.Range("E48").Value = "  -1.52%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "7.72"
$ws.Range("E49").Value = "  +0.22%  "
$ws.Range("E50").Value = "  +0.19%  "
$ws.Range("D51").Value = "1.01"
$ws.Range("E51").Value = "  +0.56%  "

Write-Host "Update complete"
